$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "74.320.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.590.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.88%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +14.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.208"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +27.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.589.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.81"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000193"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +10.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.075.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.057.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +14.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.594.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +23.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +9.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +20.60%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +14.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.720.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +15.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "499.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +17.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +16.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.66"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.19"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.19%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.85"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +12.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.67"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +20.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.08"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +11.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0811"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +13.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.521"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.53%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.73%  "
